$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '20.552.55'
$ws.Range("E2").Value = '  +1.80%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.474.75'
$ws.Range("E3").Value = '  +2.49%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9565'
$ws.Range("E5").Value = '  +4.13%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '278.15'
$ws.Range("E6").Value = '  +0.41%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3622'
$ws.Range("E7").Value = '  -1.17%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3069'
$ws.Range("E8").Value = '  -1.91%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.70'
$ws.Range("E9").Value = '  +2.10%  '

# Row 10
$ws.Range("E10").Value = '  +4.15%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06663'
$ws.Range("E11").Value = '  +2.09%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.05%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.534'
$ws.Range("E13").Value = '  +2.48%  '

# Row 14
$ws.Range("E14").Value = '  +3.21%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.198'
$ws.Range("E15").Value = '  +2.19%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9565'
$ws.Range("E16").Value = '  +2.21%  '

# Row 17
$ws.Range("E17").Value = '  +1.32%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.476.61'
$ws.Range("E18").Value = '  +2.11%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.05929'
$ws.Range("E19").Value = '  +5.17%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.13'
$ws.Range("E20").Value = '  +1.76%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.507'

# Row 22
$ws.Range("E22").Value = '  +0.40%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.18'
$ws.Range("E23").Value = '  +3.13%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.254'
$ws.Range("E24").Value = '  -0.90%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '20.558.12'
$ws.Range("E25").Value = '  +1.74%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '143.41'
$ws.Range("E26").Value = '  +5.19%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.136'
$ws.Range("E27").Value = '  -2.36%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.20'
$ws.Range("E28").Value = '  +1.45%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.637.82'
$ws.Range("E29").Value = '  +2.41%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '113.79'
$ws.Range("E30").Value = '  +2.89%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.922'
$ws.Range("E31").Value = '  +4.10%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.990'
$ws.Range("E32").Value = '  +3.39%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8107'
$ws.Range("E33").Value = '  +0.05%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07999'
$ws.Range("E34").Value = '  +4.05%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.514'
$ws.Range("E35").Value = '  +2.27%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.220'
$ws.Range("E36").Value = '  +7.37%  '

# Row 37
$ws.Range("E37").Value = '  -3.67%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.743'
$ws.Range("E38").Value = '  +0.54%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02060'
$ws.Range("E39").Value = '  +3.33%  '

# Row 40
$ws.Range("E40").Value = '  +1.95%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9570'
$ws.Range("E41").Value = '  +2.33%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1880'
$ws.Range("E42").Value = '  +2.52%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.430'
$ws.Range("E43").Value = '  +2.46%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5292'
$ws.Range("E44").Value = '  +0.93%  '

# Row 45
$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.524'
$ws.Range("E45").Value = '  +0.18%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.22'
$ws.Range("E46").Value = '  +1.93%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '118.37'
$ws.Range("E47").Value = '  -0.84%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5212'
$ws.Range("E48").Value = '  +1.28%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.819'
$ws.Range("E49").Value = '  +2.74%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06481'
$ws.Range("E50").Value = '  +2.56%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9848'
$ws.Range("E51").Value = '  -0.93%  '
